$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.862.09'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.869.79'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.77'
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5080'
$ws.Range("E7").Value = '  -1.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3665'
$ws.Range("E8").Value = '  -2.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8914'
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.881.58'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07510'
$ws.Range("E13").Value = '  -1.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.84'
$ws.Range("E14").Value = '  +5.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.223'
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008500'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.919.81'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.013'
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.120.34'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.35'
$ws.Range("E23").Value = '  -1.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.376'
$ws.Range("E24").Value = '  -1.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.36'
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.777'
$ws.Range("E26").Value = '  -3.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.86'
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.091'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.34'
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.687'
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("E31").Value = '  +1.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09135'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05062'
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7477'
$ws.Range("E34").Value = '  +2.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.982'
$ws.Range("E35").Value = '  -2.86%  '
$ws.Range("E36").Value = '  -0.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.230'
$ws.Range("E37").Value = '  +5.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.530'
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5595'
$ws.Range("E39").Value = '  +5.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01993'
$ws.Range("E40").Value = '  -2.08%  '
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.625'
$ws.Range("E42").Value = '  +2.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '115.95'
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.564'
$ws.Range("E44").Value = '  +3.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1478'
$ws.Range("E45").Value = '  +0.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4773'
$ws.Range("E46").Value = '  +2.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9999'
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.08'
$ws.Range("E48").Value = '  +1.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.559'
$ws.Range("E49").Value = '  -0.65%  '
$ws.Range("E50").Value = '  +1.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.03'
$ws.Range("E51").Value = '  -0.91%  '
